$wb = $excel.ActiveWorkbook

# --- Sheet "3.Materia-Curso": fix formulas "cursos" -> "curso" ---
$ws3 = $wb.Worksheets.Item("3.Materia-Curso")

$ws3.Range("F1").Formula = "=CONCATENATE(""INSERT into curso SET idCurso="",A1,"", nombreCurso='"",B1,`n""', anio="",C1,"", idPeriodo="",D1,"", dniProfesor='"",E1,""' ;"")"

$ws3.Range("F2").Formula = "=CONCATENATE(""INSERT into curso SET idCurso="",A2,"", nombreCurso='"",B2,`n""', anio="",C2,"", idPeriodo="",D2,"", dniProfesor='"",E2,""' ;"")"
$ws3.Range("F3").Formula = "=CONCATENATE(""INSERT into curso SET idCurso="",A3,"", nombreCurso='"",B3,`n""', anio="",C3,"", idPeriodo="",D3,"", dniProfesor='"",E3,""' ;"")"
$ws3.Range("F4").Formula = "=CONCATENATE(""INSERT into curso SET idCurso="",A4,"", nombreCurso='"",B4,`n""', anio="",C4,"", idPeriodo="",D4,"", dniProfesor='"",E4,""' ;"")"
$ws3.Range("F5").Formula = "=CONCATENATE(""INSERT into curso SET idCurso="",A5,"", nombreCurso='"",B5,`n""', anio="",C5,"", idPeriodo="",D5,"", dniProfesor='"",E5,""' ;"")"

# --- Selections / views ---
$ws1 = $wb.Worksheets.Item("1.TiposDe")
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 10
$ws1.Range("B26").Select()

$ws2 = $wb.Worksheets.Item("2.Usuarios")
$ws2.Activate()
$ws2.Range("M1:M12").Select()

$ws3.Activate()
$ws3.Range("F2").Select()

$wb.Activate()
